$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The vendor pincode mapping template is being trimmed down: the
# Vendor_ID, Appliance_ID, Brand, Area and Region columns are removed,
# leaving Vendor_Name, Appliance, Pincode, City and State.
# Shift the columns that should survive (C, G, I, J) left into
# B, C, D, E respectively (A already holds Vendor_Name / vendor_name).

$ws.Range("B1").Value = $ws.Range("C1").Value2
$ws.Range("B2").Value = $ws.Range("C2").Value2

$ws.Range("C1").Value = $ws.Range("G1").Value2
$ws.Range("C2").Value = $ws.Range("G2").Value2

$ws.Range("D1").Value = $ws.Range("I1").Value2
$ws.Range("D2").Value = $ws.Range("I2").Value2

$ws.Range("E1").Value = $ws.Range("J1").Value2
$ws.Range("E2").Value = $ws.Range("J2").Value2

# Wipe the now-unused trailing columns (F:J) entirely.
$ws.Range("F1:J2").Clear() | Out-Null

# Move the active selection like the saved workbook shows.
$ws.Range("H8").Select() | Out-Null
